$d = $word.ActiveDocument
$parts = @(
'```plaintext',
'[Your Name Here]',
'[Address] | [City, State, Zip] | [LinkedIn Profile URL]  ',
'Phone: (+65) 9326 1620  ',
'Email: qixianggoh@gmail.com  ',
'',
'------------------------------------------------------------',
'',
'**Career Summary**  ',
'Driven and analytical Business Analyst with a strong technical proficiency in data analysis tools such as Python and dashboarding software, cultivated through self-driven learning. Possesses a comprehensive understanding of hotel business operations and a proven track record in pricing strategies aimed at optimizing revenue streams. Actively seeking to leverage data analysis and machine learning techniques developed during a successful career transition from the hotel industry to a data analyst role, focusing on driving business growth.',
'',
'------------------------------------------------------------',
'',
'**Key Skills**  ',
'- **Programming Languages:** Proficient in Python; experience with SQL.  ',
'- **Data Visualization Tools:** Familiarity with Tableau, Power BI, and Microstrategy.  ',
'- **Machine Learning Knowledge:** Strong understanding of algorithms including hierarchical clustering and zero-shot classification.  ',
'- **Data Management Skills:** Skilled in data warehousing and ETL processes, adept in handling large datasets.  ',
'- **Statistical Modeling:** Competence in applying statistical models for data analysis.  ',
'- **Communication:** Excellent communication skills for presenting results to stakeholders, with a proven ability to provide valuable insights to executive teams.  ',
'- **Big Data Technologies:** Familiar with concepts and technologies related to big data.  ',
'',
'------------------------------------------------------------',
'',
'**Professional Experience**  ',
'',
'**Manager, Reporting and Analytics**  ',
'*Ministry of National Development* | Apr 2022 - Present  ',
'- Spearheaded data analysis projects utilizing Python and machine learning techniques for text classification, greatly enhancing the extraction of insights from unstructured data.  ',
'- Led the development of a comprehensive workflow solution for process approval routing and data capture, significantly optimizing operational efficiency.  ',
'- Designed and maintained Python scripts for financial dashboards that provided real-time monitoring and insights into financial operations.  ',
'',
'**Data Management Lead**  ',
'*Covid 19 Joint Ops Team* | Oct 2020 - Mar 2022  ',
'- Supervised a team managing data for over 70 quarantine facilities, ensuring accurate daily operational updates.  ',
'- Identified key performance indicators and provided actionable recommendations for process improvements.  ',
'- Created a streamlined data collection process utilizing web forms and Power Query, thereby improving data accuracy across the board.  ',
'- Developed dynamic dashboards for senior management evaluation of facility performance, achieving enhanced decision-making capabilities.  ',
'',
'**Revenue Analyst**  ',
'*The Ritz-Carlton Millennia Singapore* | Sep 2018 - Oct 2020  ',
'- Conducted market demand analyses to optimize revenue through strategic interventions.  ',
'- Monitored competitor activities and relayed strategic insights to the executive team for informed decision-making.  ',
'- Ensured timely and accurate reporting of insights, greatly facilitating decision-making processes.  ',
'',
'------------------------------------------------------------',
'',
'**Prior Experience**  ',
'- **Revenue Analyst** | Royal Plaza on Scotts & 8 On Claymore (Oct 2016 - Aug 2018)  ',
'- **Reservation Sales Executive** | The Ritz-Carlton Millennia Singapore (Jan 2015 - Sep 2016)  ',
'- **Front Desk Agent** | The Ritz-Carlton Millennia Singapore (Aug 2014 - Jan 2015)  ',
'',
'------------------------------------------------------------',
'',
'**Education**  ',
'- **Bachelor of Science in Hotel Administration** | 2012-2014  ',
'  Singapore Institute of Technology, conferred by University of Nevada, Las Vegas  ',
'- **Diploma in Sports & Exercise Sciences** | 2007-2010  ',
'  Republic Polytechnic  ',
'',
'------------------------------------------------------------',
'',
'**Professional Development**  ',
'- Ongoing self-driven learning in data analysis and machine learning.  ',
'- Actively pursuing industry certifications to enhance qualifications (e.g., Microsoft Certified: Azure Data Scientist Associate).  ',
'',
'------------------------------------------------------------',
'',
'**Interests**  ',
'- Passionate about data-driven decision-making and implementing machine learning principles within business contexts.  ',
'- Enthusiastic about exploring emerging technologies and data analysis methods that promote operational efficiency and growth.  ',
'',
'------------------------------------------------------------',
'',
'**Communication Style**  ',
'- Versatile communicator, adept at collaborating with cross-functional teams and engaging non-technical stakeholders.  ',
'- Demonstrates clarity and precision in delivering complex data narratives through dashboards and reports.  ',
'- Proven ability to present analytical findings compellingly, influencing informed decision-making processes.  ',
'',
'------------------------------------------------------------',
'```',
'',
'This updated resume is designed to effectively highlight Goh Qi Xiang''s qualifications and experiences that align perfectly with prospective data analyst or data scientist roles, ensuring it stands out in the competitive job market. The structured format, organized sections, and clear articulation of skills and experience will resonate well with hiring managers.'
)
$newText = $parts -join [char]11
$d.Paragraphs(2).Range.Text = $newText
